$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.653.54'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Value = '3.898.60'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Value = '''604.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').Value = '''171.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.52%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').Value = '3.900.47'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.44%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('E10').Value = '  +1.29%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').Value = '''6.41'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('E12').Value = '  +1.74%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').Value = '''0.0000257'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.91%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Value = '''38.35'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.59%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Value = '4.550.78'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').Value = '3.898.30'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Value = '69.704.04'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.24%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Value = '''18.68'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +9.12%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('E20').Value = '  -0.74%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').Value = '''11.07'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Value = '''491.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').Value = '''0.748'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.03%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').Value = '''0.0000165'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.39%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').Value = '''85.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.57%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').Value = '''2.31'
$ws.Range('D26').Style = 'Normal'

$ws.Range('D27').Value = '''12.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.93%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').Value = '''10.17'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D30').Value = '''2.99'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').Value = '4.046.99'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('E32').Value = '  +1.44%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').Value = '''7.87'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').Value = '''31.99'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').Value = '3.864.54'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.88%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('E36').Value = '  -0.31%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').Value = '''3.44'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +15.01%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('E38').Value = '  +4.00%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('E40').Value = '  +0.42%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('D42').Value = '''0.329'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.00%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('E43').Value = '  +5.27%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').Value = '''435.59'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').Value = '''8.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.40%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('E48').Value = '  +20.87%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('E49').Value = '  +2.61%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').Value = '''40.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.99%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').Value = '''142.99'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.15%  '
$ws.Range('E51').Style = 'Normal'
